$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.55"
$ws.Range("E2").Value = "'-0.19%"
$ws.Range("D3").Value = "'36.70"
$ws.Range("E3").Value = "'-2.46%"
$ws.Range("D4").Value = "'5.128"
$ws.Range("E4").Value = "'0.00%"
$ws.Range("D5").Value = "'0.07743"
$ws.Range("E5").Value = "'-2.08%"
$ws.Range("D6").Value = "'4.412"
$ws.Range("E6").Value = "'-0.10%"
$ws.Range("D7").Value = "'8.306"
$ws.Range("E7").Value = "'0.28%"
$ws.Range("D8").Value = "'1.862"
$ws.Range("E8").Value = "'-3.60%"
$ws.Range("E9").Value = "'0.98%"
$ws.Range("D10").Value = "'0.9223"
$ws.Range("E10").Value = "'0.14%"
$ws.Range("D11").Value = "'0.1143"
$ws.Range("E11").Value = "'-7.25%"
$ws.Range("D12").Value = "'0.1875"
$ws.Range("E12").Value = "'-2.61%"
$ws.Range("D13").Value = "'0.08808"
$ws.Range("E13").Value = "'-3.67%"
$ws.Range("D14").Value = "'0.03382"
$ws.Range("E14").Value = "'2.08%"
$ws.Range("D15").Value = "'0.09527"
$ws.Range("E15").Value = "'-0.75%"
$ws.Range("D16").Value = "'0.001389"
$ws.Range("E16").Value = "'-0.05%"
$ws.Range("D17").Value = "'0.005848"
$ws.Range("E17").Value = "'0.92%"
$ws.Range("D18").Value = "'3.357"
$ws.Range("E18").Value = "'-4.43%"
$ws.Range("D19").Value = "'0.3432"
$ws.Range("E19").Value = "'-0.37%"
$ws.Range("D20").Value = "'6.327"
$ws.Range("E20").Value = "'20.46%"
$ws.Range("D21").Value = "'0.1290"
$ws.Range("E21").Value = "'1.35%"
$ws.Range("D23").Value = "'0.04330"
$ws.Range("E23").Value = "'-0.58%"
$ws.Range("D24").Value = "'0.001207"
$ws.Range("E24").Value = "'-3.38%"
$ws.Range("D25").Value = "'0.004249"
$ws.Range("E25").Value = "'-1.27%"
$ws.Range("D26").Value = "'0.0001338"
$ws.Range("E26").Value = "'9.71%"
$ws.Range("D27").Value = "'0.0002914"
$ws.Range("E27").Value = "'-98.61%"
$ws.Range("D39").Value = "'0.02091"
$ws.Range("E39").Value = "'-7.24%"
$ws.Range("D40").Value = "'0.05035"
$ws.Range("E40").Value = "'-1.67%"
$ws.Range("D41").Value = "'0.007508"
$ws.Range("E41").Value = "'0.77%"
$ws.Range("D42").Value = "'0.1349"
$ws.Range("E42").Value = "'-1.01%"
$ws.Range("D43").Value = "'0.008481"
$ws.Range("E43").Value = "'-3.47%"
$ws.Range("D44").Value = "'0.002008"
$ws.Range("E44").Value = "'-0.04%"
$ws.Range("D45").Value = "'0.007753"
$ws.Range("E45").Value = "'-9.96%"
$ws.Range("D46").Value = "'0.00006400"
$ws.Range("E46").Value = "'-4.98%"
$ws.Range("D47").Value = "'0.00000000754"
$ws.Range("E47").Value = "'0.49%"
$ws.Range("D48").Value = "'0.002882"
$ws.Range("E48").Value = "'-13.96%"
$ws.Range("D49").Value = "'0.001697"
$ws.Range("E49").Value = "'41.43%"
$ws.Range("D50").Value = "'0.00002111"
$ws.Range("E50").Value = "'0.49%"
$ws.Range("D51").Value = "'0.0002010"
$ws.Range("E51").Value = "'0.49%"
